# Outliers Design Document.docx - apply commit edit
#
# Summary of the target change (from the OOXML diff):
#  1) "The timer is designed as a way to impose ..." loses " as a way",
#     becoming "The timer is designed to impose ...". The resulting
#     paragraph ends up split into three runs ("The timer is designed ",
#     "to", " impose ...") with the "_GoBack" bookmark sitting right
#     after "to" (this is simply Word's automatic "last edit" bookmark
#     landing where the edit happened).
#  2) In the Controls paragraph, the "_GoBack" bookmark that used to sit
#     between "should" and "provide" is gone, and two new runs appear in
#     its place: an extra space run, and a run containing a single
#     <w:tab/> character.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Edit 1: "designed as a way to impose" -> "designed to impose"
# ---------------------------------------------------------------------
$asAWay = $d.Content.Duplicate
$asAWay.Find.Execute(" as a way")
$asAWay.Delete()

# Locate "designed to impose" (post-delete) so we can split the run right
# after "designed " (9 chars) and right after "to" (11 chars total).
$designedTo = $d.Content.Duplicate
$designedTo.Find.Execute("designed to impose")

# Force a run boundary between "designed " and "to" using a throwaway
# bookmark (inserting/placing a bookmark inside a run's text splits that
# run at the bookmark's position). We remove this helper bookmark again
# once the split has been created.
$splitPos = $designedTo.Start + 9
$splitRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("TempSplit", $splitRange)

# ---------------------------------------------------------------------
# Edit 2: insert an extra space + a literal tab run right before
# "provide feedback ..." in the Controls paragraph (replacing the single
# space run that currently sits there, which also currently owns the
# "_GoBack" bookmark).
# ---------------------------------------------------------------------
$provide = $d.Content.Duplicate
$provide.Find.Execute(" provide feedback")
$spaceRange = $d.Range($provide.Start, $provide.Start + 1)

$openXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body><w:p>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
  '<w:r><w:tab/></w:r>' +
  '</w:p></w:body></w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'
$spaceRange.InsertXML($openXml)

# ---------------------------------------------------------------------
# Move the "_GoBack" bookmark so it lands right after "to" (between "to"
# and " impose"). Because a document can only have one bookmark with a
# given name, adding it here automatically removes it from wherever it
# used to be (the Controls paragraph).
# ---------------------------------------------------------------------
$goBackPos = $designedTo.Start + 11
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# Clean up the helper bookmark used to force the run split.
$d.Bookmarks("TempSplit").Delete()
